# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.34 = 17212.33 pesos`n✅ 17212.33 pesos = 4.3 = 925.17 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $text

# --- Sheet "tasas": update the N10/O10 and N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 230.3
$wsTasas.Range("O10").Value = 3964
$wsTasas.Range("N12").Value = 3999.99
$wsTasas.Range("O12").Value = 215
